$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Secant")

# Secant method uses two starting guesses in B3/B4; everything below them
# (B5:D28) recomputes automatically off these. Tighten the guesses.
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2.4

# Move the active selection from C5 to B5 to match the saved view state.
$ws.Activate()
$ws.Range("B5").Select()
